# Update cryptos list: prices and hourly volume percentages, plus a
# handful of re-ranked coin rows (B/C/D/E) per the Oct 31 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column (D) as plain text so Excel does not silently
# coerce numeric-looking strings (e.g. "1.00", "0.999") into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.501.11"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "1.809.56"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "225.70"
$ws.Range("E5").Value = "  -1.06%  "

$ws.Range("D6").Value = "0.600"
$ws.Range("E6").Value = "  +2.83%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").Value = "36.32"
$ws.Range("E8").Value = "  +3.81%  "

$ws.Range("E9").Value = "  -1.77%  "

$ws.Range("D10").Value = "0.0681"
$ws.Range("E10").Value = "  -1.61%  "

$ws.Range("D11").Value = "0.0967"
$ws.Range("E11").Value = "  +1.39%  "

$ws.Range("D12").Value = "2.069.34"
$ws.Range("E12").Value = "  +0.33%  "

$ws.Range("D13").Value = "11.34"
$ws.Range("E13").Value = "  +1.51%  "

$ws.Range("D14").Value = "1.829.06"
$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("E15").Value = "  -1.68%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.462.91"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "4.42"
$ws.Range("E17").Value = "  +1.86%  "

$ws.Range("D18").Value = "68.67"
$ws.Range("E18").Value = "  -0.43%  "

$ws.Range("D19").Value = "243.25"
$ws.Range("E19").Value = "  -0.89%  "

$ws.Range("D20").Value = "0.0₃0774"
$ws.Range("E20").Value = "  -2.80%  "

$ws.Range("D21").Value = "11.25"
$ws.Range("E21").Value = "  -2.20%  "

$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("E23").Value = "  -1.05%  "

$ws.Range("E24").Value = "  +5.29%  "

$ws.Range("D25").Value = "171.25"
$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("D26").Value = "7.91"
$ws.Range("E26").Value = "  +3.64%  "

$ws.Range("D27").Value = "17.29"
$ws.Range("E27").Value = "  +3.35%  "

$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("D31").Value = "3.93"
$ws.Range("E31").Value = "  -1.63%  "

$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("E33").Value = "  -1.98%  "

$ws.Range("E34").Value = "  -1.87%  "

$ws.Range("E35").Value = "  -2.42%  "

$ws.Range("D36").Value = "0.654"
$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("D38").Value = "2.37"

$ws.Range("E39").Value = "  -1.87%  "

$ws.Range("E40").Value = "  +1.62%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.79"
$ws.Range("E41").Value = "  -1.63%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "81.13"
$ws.Range("E42").Value = "  -2.31%  "

$ws.Range("D43").Value = "0.938"
$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").Value = "1.17"
$ws.Range("E44").Value = "  +4.91%  "

$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("E46").Value = "  -1.79%  "

$ws.Range("D47").Value = "1.971.40"
$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("E48").Value = "  -2.77%  "

$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("D50").Value = "102.75"
$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0122"
$ws.Range("E51").Value = "  -6.70%  "

# Remove the temporary text-format style so the cells end up with no
# explicit style index, matching their original (default) formatting.
$ws.Range("D2:D51").Style = "Normal"
